$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO" (product-group pivot)
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("H3").Value = 595.08
$wsGrupo.Range("M9").Value = 400.46
$wsGrupo.Range("H12").Value = 102.6
$wsGrupo.Range("H14").Value = "2 de 12"

# Sheet "VENTA MENSUAL" (monthly pivot, column F = septiembre)
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F3").Value = 595.08
$wsMensual.Range("F9").Value = 400.46
$wsMensual.Range("F12").Value = 102.6
$wsMensual.Range("F14").Value = 6606.64
